# Add three new papers (rows) to the stakeholder-management literature dump.
# Existing rows 12-13 ("Descriptive stakeholder theory" / "Normative stakeholder
# theory" in column D) get pushed down to rows 15-16 (order swapped), and five
# brand-new rows of paper notes are inserted/appended around them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 12 - Buysse & Verbeke (2003), "Pays to be green"
# ---------------------------------------------------------------------------
$ws.Cells.Item(12, 1).Value = "Buysse, K., & Verbeke, A. (2003). Proactive Environmental Strategies: A Stakeholder Management Perspective. Strategic Management Journal, 24(5), 453–470."
$ws.Cells.Item(12, 2).Value = 2003
$ws.Cells.Item(12, 3).Value = "Managing for multiple stakeholders relationship to environmental strategy"
$ws.Cells.Item(12, 4).Value = "Pays to be green"
$ws.Cells.Item(12, 6).Value = "Different stakeholders matter in different ways to firms' corporate environmetnal strategy. Data are from Belgian manufacturing firms that are mostly B2B mid-supply chain operators. Internal, primary stakeholders seemed most important to the firms that had an enviro leadership strategy."
$ws.Cells.Item(12, 7).Value = "Overall, firms perceived regulators and international agreements as the most important stakeholder groups for their corporate environmetnal policy formulation."
$ws.Cells.Item(12, 8).Value = "Firms with enviro leadership strategies seemed to care only about internal primary stakeholders and not other stakeholder groups. This could result from the studied firms (Belgian manufacturing firms) being B2B rather than B2C so they don't interact with external primary stakeholders like customers."
$ws.Rows.Item(12).RowHeight = 75

# ---------------------------------------------------------------------------
# Row 13 - Tashman & Raelin (2013), "Stakeholder-agency theory"
# ---------------------------------------------------------------------------
$ws.Cells.Item(13, 1).Value = "Tashman, P., & Raelin, J. (2013). Who and What Really Matters to the Firm: Moving Stakeholder Salience beyond Managerial Perceptions. Business Ethics Quarterly, 23(4), 591–616. http://doi.org/10.5840/beq201323441"
$ws.Cells.Item(13, 2).Value = 2013
$ws.Cells.Item(13, 3).Value = "Stakeholder salience to the firm"
$ws.Cells.Item(13, 4).Value = "Stakeholder-agency theory"
$ws.Cells.Item(13, 5).Value = "The concept of stakeholder salience is based on managerial perception, but some stakeholders should matter to the firm even when managers don't perceive them as important. They develop the concept of stakeholder salience to the firm that expands salience evaluation beyond firm managers to include the societal level of analysis, such that stakeholder salience is a function of both organization-level and society-level perceptions."
$ws.Cells.Item(13, 6).Value = "Move salience from salience to managers [salience = f (manager perception)] to salience to the firm [salience = f (manager perception, stakeholder perception)]. Managers that correctly assess stakeholder salience should have improved firm performance, but empirical studies find mixed results. They argue this is because market frictions can prevent managers from correctly assessing stakeholder salience. Stakeholder dialogue can reduce the negative affect of frictions."
$ws.Cells.Item(13, 7).Value = "Models testing stakeholder salience are misspecified if they include only managers' assessment of salience and omit stakeholders' assessment of salience of themselves and other stakeholders."
$ws.Cells.Item(13, 8).Value = "Market frictions (they draw on nexus of contracts here so market frictions are problems that prevent complete contracting) that are low should make manager assessment better, such that market frictions positively moderate the manager-based salience firm performance link. When market frictions are high, they negatively moderate the link. But engaging in stakeholder dialogue is a moderator of the moderation effect and can reduce the negative moderation of high market frictions."
$ws.Cells.Item(13, 9).Value = "Mitchell, Agle, and Wood's (1997) original formualtion of salience around power, urgency, and legitimacy needs a stronger normative component to prevent powerful stakeholders from dominating managers' attention. The refined concept of stakeholder salience to the firm includes a normative component by recognizing hypernorms, institutional expectations, and stakeholder perceptions of legitimacy and urgency."
$ws.Cells.Item(13, 10).Value = "How to assess perceptions of salience across multiple stakeholder groups. How to modify scales developed to assess managers' perceptions of power, legitimacy, and urgency to assess other stakeholder groups' perceptions of these characteristics?"
$ws.Rows.Item(13).RowHeight = 135

# ---------------------------------------------------------------------------
# Row 14 - Convergent stakeholder theory (continuation of Tashman & Raelin)
# ---------------------------------------------------------------------------
$ws.Cells.Item(14, 4).Value = "Convergent stakeholder theory"
$ws.Cells.Item(14, 6).Value = "`"The fundamental goal of convergent stakeholder theory is to reconcile its normative and instrumental strands and prescribe decision making that helps firms prosper while treating stakeholders ethically`" (Tashman & Raelin (2013)"
$ws.Rows.Item(14).RowHeight = 45

# ---------------------------------------------------------------------------
# Row 15 / 16 - the previous rows 13/12 (Normative / Descriptive stakeholder
# theory), now pushed down and re-ordered.
# ---------------------------------------------------------------------------
$ws.Cells.Item(15, 4).Value = "Normative stakeholder theory"
$ws.Cells.Item(16, 4).Value = "Descriptive stakeholder theory"

# ---------------------------------------------------------------------------
# Row 17 - Hall, Millo & Barman (2015), "Stakeholder management"
# ---------------------------------------------------------------------------
$ws.Cells.Item(17, 1).Value = "Hall, M., Millo, Y., & Barman, E. (2015). Who and What Really Counts? Stakeholder Prioritization and Accounting for Social Value. Journal of Management Studies, 52(7), 907–934. http://doi.org/10.1111/joms.12146"
$ws.Cells.Item(17, 2).Value = 2015
$ws.Cells.Item(17, 3).Value = "Stakeholder management"
$ws.Cells.Item(17, 4).Value = "Stakeholder theory"
$ws.Cells.Item(17, 5).Value = "Stakeholder management research has not addressed how organizational practices influence stakeholder prioritization."

$f17 = $ws.Cells.Item(17, 6)
$f17.Value = "A firm's accounting and reporting system influences which stakeholders view as salient to the firm, demonstrating that stakeholder prioritization is not driven solely by managerial decisions but is influenced by the systems used in the process of stakeholder prioritization decisions."
$f17.Characters(120, 164).Font.Bold = $true

$ws.Cells.Item(17, 7).Value = "Two research questions: (1) How do managers develop an accounting and reporting system to reflect their prioritizaiton of stakeholders? (2) What factors influence managers' ability to construct an accounting and reporting system to incorporate the voices of salient stakeholders?"
$ws.Cells.Item(17, 8).Value = "They study the Social Return on Investemnt accounting methodology and compare cases in the UK and US. They find that organizations customize the SROI system in ways that reflects different assumptions and resources available to managers."
$ws.Cells.Item(17, 9).Value = "From their findings, they enrich stakeholder prioritization theory by noting the importance of managers' assumptions about which types of knowledge are valid and acceptable and of the technical and matrial resources available to managers to assess stakeholders."
$ws.Cells.Item(17, 10).Value = "Test the two propositions they develop: (1) Prioritization of stakeholder voices (which stakeholders' voices are included, how are they represented and measured) is influenced by managers' epistemic beliefs (such as what counts as valid and appropriate data). (2) Ability of managers to develop accounting and reporting system consisten with their epistemic beliefs is shaped by organization's material conditions (nature of existing data collection and reporting systems, access to financial resources, access to necessary labor and expertise)."
$ws.Rows.Item(17).RowHeight = 120

# ---------------------------------------------------------------------------
# Column width tweaks: A narrows, G widens (to fit the new notes).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 40.7
$ws.Columns.Item(7).ColumnWidth = 58.2

# ---------------------------------------------------------------------------
# Move the view: frozen pane scrolled down to row 11 / col J, selection on J17.
# ---------------------------------------------------------------------------
$ws.Range("J17").Select()
